# Add "Add Two Numbers" problem, restructure the table:
#  - drop the "Notes" and "URL" text columns
#  - add a single "Local Path" hyperlink column instead
#  - move "Revisit" to sit right after "Blind75 problem?"
#  - append the new "Add Two Numbers" / Linked List row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the old "URL" column (H) entirely; "Notes" (G) gets
#    overwritten below with the new "Local Path" hyperlink column.
# ------------------------------------------------------------------
$ws.Columns.Item(8).Delete()

# ------------------------------------------------------------------
# 2. Rewrite the whole table in the new column order:
#    A Name | B Category | C Blind75? | D Revisit | E Difficulty |
#    F Relative Difficulty | G Local Path
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Blind75 problem?"
$ws.Range("D1").Value = "Revisit"
$ws.Range("E1").Value = "Difficulty"
$ws.Range("F1").Value = "Relative Difficulty"
$ws.Range("G1").Value = "Local Path"

$ws.Range("A2").Value = "TwoSum"
$ws.Range("B2").Value = "Array"
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "Easy"
$ws.Range("F2").Value = "Easy"
$ws.Range("G2").Value = "1 - Two Sum"

$ws.Range("A3").Value = "Longest Common Prefix"
$ws.Range("B3").Value = "String"
$ws.Range("C3").Value = "No"
$ws.Range("D3").Value = "No"
$ws.Range("E3").Value = "Easy"
$ws.Range("F3").Value = "Easy"
$ws.Range("G3").Value = "14 - Longest Common Prefix"

$ws.Range("A4").Value = "Valid Parentheses"
$ws.Range("B4").Value = "String"
$ws.Range("C4").Value = "Yes"
$ws.Range("D4").Value = "No"
$ws.Range("E4").Value = "Easy"
$ws.Range("F4").Value = "Easy"
$ws.Range("G4").Value = "20 - Valid Parentheses"

$ws.Range("A5").Value = "Add Two Numbers"
$ws.Range("B5").Value = "Linked List"
$ws.Range("C5").Value = "No"
$ws.Range("D5").Value = "No"
$ws.Range("E5").Value = "Easy"
$ws.Range("F5").Value = "Easy"
$ws.Range("G5").Value = "2 - Add Two Numbers"

# ------------------------------------------------------------------
# 3. Hyperlinks for the "Local Path" column.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G2"), "1 - Two Sum")
$ws.Hyperlinks.Add($ws.Range("G3"), "14 - Longest Common Prefix")
$ws.Hyperlinks.Add($ws.Range("G4"), "20 - Valid Parentheses")
$ws.Hyperlinks.Add($ws.Range("G5"), "2 - Add Two Numbers")

# ------------------------------------------------------------------
# 4. Column widths for the new A:G layout.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24.31   # A Name
$ws.Columns.Item(2).ColumnWidth = 8.74    # B Category
$ws.Columns.Item(3).ColumnWidth = 17.02   # C Blind75 problem?
$ws.Columns.Item(4).ColumnWidth = 15.88   # D Revisit
$ws.Columns.Item(5).ColumnWidth = 17.59   # E Difficulty
$ws.Columns.Item(6).ColumnWidth = 18.02   # F Relative Difficulty
$ws.Columns.Item(7).ColumnWidth = 30.31   # G Local Path

# ------------------------------------------------------------------
# 5. Conditional formatting + data validation ranges need to grow to
#    cover the new row 5 (and, for validation, the shifted columns).
# ------------------------------------------------------------------
$fcs = $ws.Range("D2:G4").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("D2:G5"))
}

$ws.Cells.Validation.Delete()
$ws.Range("E2:F5").Validation.Add(3, 1, 0, '"Easy, Medium, Hard"')
$ws.Range("C2:C5").Validation.Add(3, 1, 0, '"Yes, No"')
$ws.Range("C2:C5").Validation.IgnoreBlank = $false
$ws.Range("B2:B5").Validation.Add(3, 1, 0, '"Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap"')
$ws.Range("D2:D5").Validation.Add(3, 1, 0, '"Yes, No"')

# ------------------------------------------------------------------
# 6. Selection, matching where the author's cursor ended up.
# ------------------------------------------------------------------
$ws.Range("G9").Select() | Out-Null
